# Generate Report for Handoff
# The file e2e\75af15d5-9871-42fc-9627-421b87f9cd98.md (row 3 in every sheet)
# moves from "Handed back: in sync with en-US" to "Ready for handoff" status,
# gets a fresh handoff datetime, and picks up an error detail message about
# a stale handback file.

$wb = $excel.ActiveWorkbook

$status = "Ready for handoff"
$overviewDate = "2016-08-18 20:49:42"
$zhHandoffDate = "2016-08-18 20:49:37"
$deHandoffDate = "2016-08-18 20:49:42"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/eb3b58a3e98790a53ba9044c088470a3af38b2fa/e2e/75af15d5-9871-42fc-9627-421b87f9cd98.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7f2a1bb25135b936f0aba60848438928d96e0997/e2e/75af15d5-9871-42fc-9627-421b87f9cd98.md."

# --- Overview sheet: update the status columns for the second file ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $status
$overview.Range("F3").Value = $status
$overview.Range("G3").Value = $overviewDate

# --- zh-cn detail sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $status
$zhcn.Range("H3").Value = $zhHandoffDate
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Columns.Item(16).ColumnWidth = 39.16666666666667

# --- de-de detail sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $status
$dede.Range("H3").Value = $deHandoffDate
$dede.Range("P3").Value = $errorDetail
$dede.Columns.Item(16).ColumnWidth = 39.16666666666667
